# Auto-generated edit script applying scheduled-runner price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 589.8
$ws.Range("I55").Value = 672.5
$ws.Range("J55").Value = 259
$ws.Range("K55").Value = 672.5
$ws.Range("L55").Value = 259
$ws.Range("M55").Value = -458.5
$ws.Range("N55").Value = -687
$ws.Range("H64").Value = 111118056
$ws.Range("J64").Value = 125006690
$ws.Range("L64").Value = 125006690
$ws.Range("N64").Value = -125007186
$ws.Range("H67").Value = 111118056
$ws.Range("J67").Value = 125006690
$ws.Range("L67").Value = 125006690
$ws.Range("N67").Value = -125008406
$ws.Range("H86").Value = 58826470
$ws.Range("J86").Value = 4035.7144
$ws.Range("L86").Value = 4035.7144
$ws.Range("N86").Value = -6281.7144
$ws.Range("H88").Value = 9091859
$ws.Range("J88").Value = 1068.75
$ws.Range("L88").Value = 1068.75
$ws.Range("N88").Value = -1880.75
$ws.Range("H89").Value = 58826470
$ws.Range("J89").Value = 4035.7144
$ws.Range("L89").Value = 20178.572
$ws.Range("N89").Value = -31410.572
$ws.Range("H91").Value = 9091859
$ws.Range("J91").Value = 1068.75
$ws.Range("L91").Value = 1068.75
$ws.Range("N91").Value = -3876.75
$ws.Range("H100").Value = 2311
$ws.Range("I100").Value = 1350
$ws.Range("K100").Value = 1350
$ws.Range("M100").Value = -809
$ws.Range("H112").Value = 2102.0667
$ws.Range("J112").Value = 2102.0667
$ws.Range("L112").Value = 6306.2001
$ws.Range("N112").Value = -8522.2001
$ws.Range("H113").Value = 3708.3635
$ws.Range("J113").Value = 3888.111
$ws.Range("L113").Value = 3888.111
$ws.Range("N113").Value = -10396.111
$ws.Range("H137").Value = 2278782.5
$ws.Range("I137").Value = 2941953.8
$ws.Range("J137").Value = 24000.6
$ws.Range("K137").Value = 8825861.399999999
$ws.Range("L137").Value = 72001.79999999999
$ws.Range("M137").Value = -8823311.399999999
$ws.Range("N137").Value = -77101.79999999999
$ws.Range("H138").Value = 3390.7593
$ws.Range("I138").Value = 4281.7646
$ws.Range("J138").Value = 2981.3784
$ws.Range("K138").Value = 12845.2938
$ws.Range("L138").Value = 8944.135200000001
$ws.Range("M138").Value = -7705.293800000001
$ws.Range("N138").Value = -19224.1352

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 752.6818
$ws.Range("I2").Value = 415.84616
$ws.Range("J2").Value = 1239.2222
$ws.Range("K2").Value = 415.84616
$ws.Range("L2").Value = 1239.2222
$ws.Range("M2").Value = -302.84616
$ws.Range("N2").Value = -1465.2222
$ws.Range("H32").Value = 3775.2744
$ws.Range("I32").Value = 2883.9512
$ws.Range("K32").Value = 2883.9512
$ws.Range("M32").Value = -2596.9512
$ws.Range("H61").Value = 2413.2222
$ws.Range("I61").Value = 1778.7
$ws.Range("J61").Value = 5585.8335
$ws.Range("K61").Value = 1778.7
$ws.Range("L61").Value = 5585.8335
$ws.Range("M61").Value = -1566.7
$ws.Range("N61").Value = -6009.8335
$ws.Range("H74").Value = 254843.64
$ws.Range("I74").Value = 398230.72
$ws.Range("K74").Value = 398230.72
$ws.Range("M74").Value = -397356.72
$ws.Range("H77").Value = 254843.64
$ws.Range("I77").Value = 398230.72
$ws.Range("K77").Value = 1991153.6
$ws.Range("M77").Value = -1986785.6
$ws.Range("H116").Value = 752.6818
$ws.Range("I116").Value = 415.84616
$ws.Range("J116").Value = 1239.2222
$ws.Range("K116").Value = 415.84616
$ws.Range("L116").Value = 1239.2222
$ws.Range("M116").Value = 1878.15384
$ws.Range("N116").Value = -5827.2222
$ws.Range("H132").Value = 1848.1464
$ws.Range("I132").Value = 1009.19354
$ws.Range("J132").Value = 4448.9
$ws.Range("K132").Value = 3027.58062
$ws.Range("L132").Value = 13346.7
$ws.Range("M132").Value = -497.5806199999997
$ws.Range("N132").Value = -18406.7
$ws.Range("H136").Value = 2413.2222
$ws.Range("I136").Value = 1778.7
$ws.Range("J136").Value = 5585.8335
$ws.Range("K136").Value = 5336.1
$ws.Range("L136").Value = 16757.5005
$ws.Range("M136").Value = -2786.1
$ws.Range("N136").Value = -21857.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 752.6818
$ws.Range("I3").Value = 415.84616
$ws.Range("J3").Value = 1239.2222
$ws.Range("K3").Value = 415.84616
$ws.Range("L3").Value = 1239.2222
$ws.Range("M3").Value = -301.84616
$ws.Range("N3").Value = -1467.2222
$ws.Range("H105").Value = 8968087
$ws.Range("I105").Value = 528876.1
$ws.Range("J105").Value = 25002588
$ws.Range("K105").Value = 528876.1
$ws.Range("L105").Value = 25002588
$ws.Range("M105").Value = -527129.1
$ws.Range("N105").Value = -25006082

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3347.0364
$ws.Range("I31").Value = 2026.7858
$ws.Range("J31").Value = 7612.4614
$ws.Range("K31").Value = 2026.7858
$ws.Range("L31").Value = 7612.4614
$ws.Range("M31").Value = -1731.7858
$ws.Range("N31").Value = -8202.4614
$ws.Range("H34").Value = 3347.0364
$ws.Range("I34").Value = 2026.7858
$ws.Range("J34").Value = 7612.4614
$ws.Range("K34").Value = 2026.7858
$ws.Range("L34").Value = 7612.4614
$ws.Range("M34").Value = -1824.7858
$ws.Range("N34").Value = -8016.4614
$ws.Range("H35").Value = 762.25
$ws.Range("I35").Value = 762.25
$ws.Range("K35").Value = 762.25
$ws.Range("M35").Value = -468.25
$ws.Range("H58").Value = 2964.1482
$ws.Range("I58").Value = 2103.0715
$ws.Range("J58").Value = 3891.4614
$ws.Range("K58").Value = 2103.0715
$ws.Range("L58").Value = 3891.4614
$ws.Range("M58").Value = -1900.0715
$ws.Range("N58").Value = -4297.4614
$ws.Range("H62").Value = 12507888
$ws.Range("I62").Value = 14293614
$ws.Range("K62").Value = 14293614
$ws.Range("M62").Value = -14292990
$ws.Range("H65").Value = 12507888
$ws.Range("I65").Value = 14293614
$ws.Range("K65").Value = 71468070
$ws.Range("M65").Value = -71464950
$ws.Range("H105").Value = 2082.7222
$ws.Range("I105").Value = 878
$ws.Range("J105").Value = 2849.3635
$ws.Range("K105").Value = 878
$ws.Range("L105").Value = 2849.3635
$ws.Range("M105").Value = 869
$ws.Range("N105").Value = -6343.363499999999
$ws.Range("H136").Value = 2964.1482
$ws.Range("I136").Value = 2103.0715
$ws.Range("J136").Value = 3891.4614
$ws.Range("K136").Value = 6309.2145
$ws.Range("L136").Value = 11674.3842
$ws.Range("M136").Value = -3759.2145
$ws.Range("N136").Value = -16774.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2833.5
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 2833.5
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 8500.5
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -11496.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5712.2144
$ws.Range("I113").Value = 4997.25
$ws.Range("K113").Value = 4997.25
$ws.Range("M113").Value = -2827.25
$ws.Range("H132").Value = 2729.7827
$ws.Range("I132").Value = 2323.0952
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 6969.285600000001
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -4439.285600000001
$ws.Range("N132").Value = -26060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3065.6667
$ws.Range("I7").Value = 2599
$ws.Range("J7").Value = 3999
$ws.Range("K7").Value = 2599
$ws.Range("L7").Value = 3999
$ws.Range("M7").Value = -2487
$ws.Range("N7").Value = -4223
$ws.Range("H18").Value = 15984.5
$ws.Range("I18").Value = 6969
$ws.Range("K18").Value = 6969
$ws.Range("M18").Value = -6797
$ws.Range("H22").Value = 775.3333
$ws.Range("I22").Value = 717.3333
$ws.Range("J22").Value = 833.3333
$ws.Range("K22").Value = 717.3333
$ws.Range("L22").Value = 833.3333
$ws.Range("M22").Value = -422.3333
$ws.Range("N22").Value = -1423.3333
$ws.Range("H27").Value = 775.3333
$ws.Range("I27").Value = 717.3333
$ws.Range("J27").Value = 833.3333
$ws.Range("K27").Value = 717.3333
$ws.Range("L27").Value = 833.3333
$ws.Range("M27").Value = -610.3333
$ws.Range("N27").Value = -1047.3333
$ws.Range("H40").Value = 7901.6
$ws.Range("I40").Value = 7246.3335
$ws.Range("K40").Value = 7246.3335
$ws.Range("M40").Value = -7110.3335
$ws.Range("H126").Value = 3065.6667
$ws.Range("I126").Value = 2599
$ws.Range("J126").Value = 3999
$ws.Range("K126").Value = 7797
$ws.Range("L126").Value = 11997
$ws.Range("M126").Value = -5327
$ws.Range("N126").Value = -16937
$ws.Range("H136").Value = 6034
$ws.Range("I136").Value = 6053.5713
$ws.Range("J136").Value = 5999.75
$ws.Range("K136").Value = 18160.7139
$ws.Range("L136").Value = 17999.25
$ws.Range("M136").Value = -15610.7139
$ws.Range("N136").Value = -23099.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3216.4736
$ws.Range("I132").Value = 3256.5757
$ws.Range("J132").Value = 2951.8
$ws.Range("K132").Value = 9769.7271
$ws.Range("L132").Value = 8855.400000000001
$ws.Range("M132").Value = -7239.7271
$ws.Range("N132").Value = -13915.4
